$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 with Balaraju vankala's data
$ws.Range("B2").Value = "Balaraju vankala"
$ws.Range("C2").Value = 31780.8

# Delete rows 3 through 5 (Priyanka, pattabhi, and the old Balaraju row that is now redundant)
$ws.Range("A3:D5").Delete()
